$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format column A for the new rows as Text so date-like strings
# ("2024-09-02", etc.) are stored literally instead of being parsed as dates.
$ws.Range("A677:A696").NumberFormat = "@"

# Row 677
$ws.Cells.Item(677,1).Value = "2024-09-02"
$ws.Cells.Item(677,3).Value = 1964.5
$ws.Cells.Item(677,4).Value = 1806.650024414062
$ws.Cells.Item(677,5).Value = 1050.949951171875
$ws.Cells.Item(677,6).Value = 1766.300048828125
$ws.Cells.Item(677,7).Value = 1749.5
$ws.Cells.Item(677,8).Value = 38025.74987792969
$ws.Cells.Item(677,9).Value = 0
$ws.Cells.Item(677,10).Value = 173.1139691570218

# Row 678
$ws.Cells.Item(678,1).Value = "2024-09-03"
$ws.Cells.Item(678,3).Value = 1941.25
$ws.Cells.Item(678,4).Value = 1790.449951171875
$ws.Cells.Item(678,5).Value = 1068.800048828125
$ws.Cells.Item(678,6).Value = 1769.300048828125
$ws.Cells.Item(678,7).Value = 1718.75
$ws.Cells.Item(678,8).Value = 37895.55029296875
$ws.Cells.Item(678,9).Value = -0.003423984678248408
$ws.Cells.Item(678,10).Value = 172.5212295790374

# Row 679
$ws.Cells.Item(679,1).Value = "2024-09-04"
$ws.Cells.Item(679,3).Value = 1922.449951171875
$ws.Cells.Item(679,4).Value = 1785.25
$ws.Cells.Item(679,5).Value = 1056.199951171875
$ws.Cells.Item(679,6).Value = 1749.699951171875
$ws.Cells.Item(679,7).Value = 1729.550048828125
$ws.Cells.Item(679,8).Value = 37618.74926757812
$ws.Cells.Item(679,9).Value = -0.007304314708473397
$ws.Cells.Item(679,10).Value = 171.2610802242993

# Row 680
$ws.Cells.Item(680,1).Value = "2024-09-05"
$ws.Cells.Item(680,3).Value = 1933.150024414062
$ws.Cells.Item(680,4).Value = 1790.550048828125
$ws.Cells.Item(680,5).Value = 1074.900024414062
$ws.Cells.Item(680,6).Value = 1722.900024414062
$ws.Cells.Item(680,7).Value = 1720.75
$ws.Cells.Item(680,8).Value = 37761.05065917969
$ws.Cells.Item(680,9).Value = 0.003782725220059497
$ws.Cells.Item(680,10).Value = 171.9089138316784

# Row 681
$ws.Cells.Item(681,1).Value = "2024-09-06"
$ws.Cells.Item(681,3).Value = 1901.849975585938
$ws.Cells.Item(681,4).Value = 1756.099975585938
$ws.Cells.Item(681,5).Value = 1112.650024414062
$ws.Cells.Item(681,6).Value = 1730.300048828125
$ws.Cells.Item(681,7).Value = 1715
$ws.Cells.Item(681,8).Value = 37713.70007324219
$ws.Cells.Item(681,9).Value = -0.001253953084220899
$ws.Cells.Item(681,10).Value = 171.6933481189741

# Row 682
$ws.Cells.Item(682,1).Value = "2024-09-09"
$ws.Cells.Item(682,3).Value = 1894.650024414062
$ws.Cells.Item(682,4).Value = 1746.75
$ws.Cells.Item(682,5).Value = 1077.550048828125
$ws.Cells.Item(682,6).Value = 1750.400024414062
$ws.Cells.Item(682,7).Value = 1741.199951171875
$ws.Cells.Item(682,8).Value = 37455.65051269531
$ws.Cells.Item(682,9).Value = -0.006842329446480399
$ws.Cells.Item(682,10).Value = 170.5185656673748

# Row 683
$ws.Cells.Item(683,1).Value = "2024-09-10"
$ws.Cells.Item(683,3).Value = 1912.300048828125
$ws.Cells.Item(683,4).Value = 1779.099975585938
$ws.Cells.Item(683,5).Value = 1091
$ws.Cells.Item(683,6).Value = 1756.349975585938
$ws.Cells.Item(683,7).Value = 1745.150024414062
$ws.Cells.Item(683,8).Value = 37856.65014648438
$ws.Cells.Item(683,9).Value = 0.01070598503296977
$ws.Cells.Item(683,10).Value = 172.3441348792532

# Row 684
$ws.Cells.Item(684,1).Value = "2024-09-11"
$ws.Cells.Item(684,3).Value = 1910.150024414062
$ws.Cells.Item(684,4).Value = 1778.75
$ws.Cells.Item(684,5).Value = 1077.849975585938
$ws.Cells.Item(684,6).Value = 1789.349975585938
$ws.Cells.Item(684,7).Value = 1782.650024414062
$ws.Cells.Item(684,8).Value = 37910.79992675781
$ws.Cells.Item(684,9).Value = 0.001430390170918655
$ws.Cells.Item(684,10).Value = 172.5906542358

# Row 685
$ws.Cells.Item(685,1).Value = "2024-09-12"
$ws.Cells.Item(685,3).Value = 1950.449951171875
$ws.Cells.Item(685,4).Value = 1807.599975585938
$ws.Cells.Item(685,5).Value = 1083.75
$ws.Cells.Item(685,6).Value = 1838.050048828125
$ws.Cells.Item(685,7).Value = 1812.75
$ws.Cells.Item(685,8).Value = 38550.34973144531
$ws.Cells.Item(685,9).Value = 0.01686985782212682
$ws.Cells.Item(685,10).Value = 175.5022340341858

# Row 686
$ws.Cells.Item(686,1).Value = "2024-09-13"
$ws.Cells.Item(686,3).Value = 1944.099975585938
$ws.Cells.Item(686,4).Value = 1812.800048828125
$ws.Cells.Item(686,5).Value = 1089.699951171875
$ws.Cells.Item(686,6).Value = 1826.050048828125
$ws.Cells.Item(686,7).Value = 1814.099975585938
$ws.Cells.Item(686,8).Value = 38552.5498046875
$ws.Cells.Item(686,9).Value = 0.00005707012407186833
$ws.Cells.Item(686,10).Value = 175.512249968457

# Row 687
$ws.Cells.Item(687,1).Value = "2024-09-16"
$ws.Cells.Item(687,3).Value = 1950.25
$ws.Cells.Item(687,4).Value = 1811.849975585938
$ws.Cells.Item(687,5).Value = 1094.650024414062
$ws.Cells.Item(687,6).Value = 1757.849975585938
$ws.Cells.Item(687,7).Value = 1797.199951171875
$ws.Cells.Item(687,8).Value = 38385.89990234375
$ws.Cells.Item(687,9).Value = -0.004322668751821118
$ws.Cells.Item(687,10).Value = 174.7535686499565

# Row 688
$ws.Cells.Item(688,1).Value = "2024-09-17"
$ws.Cells.Item(688,3).Value = 1952.550048828125
$ws.Cells.Item(688,4).Value = 1813.75
$ws.Cells.Item(688,5).Value = 1080.300048828125
$ws.Cells.Item(688,6).Value = 1741.150024414062
$ws.Cells.Item(688,7).Value = 1848.5
$ws.Cells.Item(688,8).Value = 38346.90075683594
$ws.Cells.Item(688,9).Value = -0.001015975803798501
$ws.Cells.Item(688,10).Value = 174.5760232525807

# Row 689
$ws.Cells.Item(689,1).Value = "2024-09-18"
$ws.Cells.Item(689,3).Value = 1892.150024414062
$ws.Cells.Item(689,4).Value = 1756.5
$ws.Cells.Item(689,5).Value = 1065.800048828125
$ws.Cells.Item(689,6).Value = 1727.25
$ws.Cells.Item(689,7).Value = 1805.599975585938
$ws.Cells.Item(689,8).Value = 37454.75048828125
$ws.Cells.Item(689,9).Value = -0.023265250931541
$ws.Cells.Item(689,10).Value = 170.5144682649789

# Row 690
$ws.Cells.Item(690,1).Value = "2024-09-19"
$ws.Cells.Item(690,3).Value = 1894.199951171875
$ws.Cells.Item(690,4).Value = 1736.5
$ws.Cells.Item(690,5).Value = 1060.75
$ws.Cells.Item(690,6).Value = 1676.449951171875
$ws.Cells.Item(690,7).Value = 1877.449951171875
$ws.Cells.Item(690,8).Value = 37317.94946289062
$ws.Cells.Item(690,9).Value = -0.003652434567236724
$ws.Cells.Item(690,10).Value = 169.8916753268739

# Row 691
$ws.Cells.Item(691,1).Value = "2024-09-20"
$ws.Cells.Item(691,3).Value = 1905.75
$ws.Cells.Item(691,4).Value = 1760.050048828125
$ws.Cells.Item(691,5).Value = 1114.699951171875
$ws.Cells.Item(691,6).Value = 1662
$ws.Cells.Item(691,7).Value = 1931.449951171875
$ws.Cells.Item(691,8).Value = 38001.24975585938
$ws.Cells.Item(691,9).Value = 0.01831023147850691
$ws.Cells.Item(691,10).Value = 173.0024312283803

# Row 692
$ws.Cells.Item(692,1).Value = "2024-09-23"
$ws.Cells.Item(692,3).Value = 1896.449951171875
$ws.Cells.Item(692,4).Value = 1752.800048828125
$ws.Cells.Item(692,5).Value = 1106.699951171875
$ws.Cells.Item(692,6).Value = 1692.900024414062
$ws.Cells.Item(692,7).Value = 1920.400024414062
$ws.Cells.Item(692,8).Value = 37915.79968261719
$ws.Cells.Item(692,9).Value = -0.002248612184892999
$ws.Cells.Item(692,10).Value = 172.6134158535041

# Row 693
$ws.Cells.Item(693,1).Value = "2024-09-24"
$ws.Cells.Item(693,3).Value = 1898.599975585938
$ws.Cells.Item(693,4).Value = 1775.599975585938
$ws.Cells.Item(693,5).Value = 1098.5
$ws.Cells.Item(693,6).Value = 1660.900024414062
$ws.Cells.Item(693,7).Value = 1838.75
$ws.Cells.Item(693,8).Value = 37717.7998046875
$ws.Cells.Item(693,9).Value = -0.00522209420840627
$ws.Cells.Item(693,10).Value = 171.7120123342822

# Row 694
$ws.Cells.Item(694,1).Value = "2024-09-25"
$ws.Cells.Item(694,3).Value = 1895.300048828125
$ws.Cells.Item(694,4).Value = 1782.400024414062
$ws.Cells.Item(694,5).Value = 1088.599975585938
$ws.Cells.Item(694,6).Value = 1654.75
$ws.Cells.Item(694,7).Value = 1722.050048828125
$ws.Cells.Item(694,8).Value = 37400.95031738281
$ws.Cells.Item(694,9).Value = -0.008400529430280024
$ws.Cells.Item(694,10).Value = 170.2695405211355

# Row 695
$ws.Cells.Item(695,1).Value = "2024-09-26"
$ws.Cells.Item(695,3).Value = 1900.25
$ws.Cells.Item(695,4).Value = 1783.849975585938
$ws.Cells.Item(695,5).Value = 1094.949951171875
$ws.Cells.Item(695,6).Value = 1682.449951171875
$ws.Cells.Item(695,7).Value = 1685.900024414062
$ws.Cells.Item(695,8).Value = 37499.49938964844
$ws.Cells.Item(695,9).Value = 0.002634934979708856
$ws.Cells.Item(695,10).Value = 170.7181896894336

# Row 696
$ws.Cells.Item(696,1).Value = "2024-09-27"
$ws.Cells.Item(696,3).Value = 1906.75
$ws.Cells.Item(696,4).Value = 1808.400024414062
$ws.Cells.Item(696,5).Value = 1121
$ws.Cells.Item(696,6).Value = 1671
$ws.Cells.Item(696,7).Value = 1640.800048828125
$ws.Cells.Item(696,8).Value = 37745.10021972656
$ws.Cells.Item(696,9).Value = 0.006549442901254356
$ws.Cells.Item(696,10).Value = 171.83629872501
